# Apply updated "dSF" (column F) values as part of a data repull / push
# and mean-calculation refresh, per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F ("dSF")
$updates = @{
    4  = -8
    7  = 4
    10 = -9
    11 = 1
    12 = -2
    14 = 0
    19 = 6
    21 = 1
    29 = 1
    31 = -2
    33 = -5
    34 = 3
    35 = -3
    37 = 9
    38 = -1
    42 = -1
    43 = 1
    47 = 1
    51 = 10
    52 = 6
    53 = 9
    54 = -4
    56 = 5
    60 = 1
    61 = -1
    63 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
